# Updates the worker account-statement (Estado de Cuenta) data table.
# Previous periods for RICARDO ALVENIZ MADERO HERRERA are removed and the
# table is rebuilt in descending chronological order (2302 -> 1607), with
# the newest period (2302) carrying an adjusted "Valor Mora" of 31200.
# Two new rows are appended for ADRIANA MILENA BARCASNEGRAS ORTIZ covering
# periods 1909 and 1908.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(
    @{Row=16; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="2302"; F=31200; G=900000},
    @{Row=17; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="2301"; F=36000; G=900000},
    @{Row=18; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="2212"; F=36000; G=900000},
    @{Row=19; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="2211"; F=36000; G=900000},
    @{Row=20; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="2210"; F=36000; G=900000},
    @{Row=21; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="2209"; F=36000; G=900000},
    @{Row=22; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="2208"; F=36000; G=900000},
    @{Row=23; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="2207"; F=36000; G=900000},
    @{Row=24; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="2206"; F=36000; G=900000},
    @{Row=25; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="2205"; F=36000; G=900000},
    @{Row=26; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="2204"; F=36000; G=900000},
    @{Row=27; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="2203"; F=36000; G=900000},
    @{Row=28; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="2202"; F=36000; G=900000},
    @{Row=29; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="2201"; F=36000; G=900000},
    @{Row=30; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="2112"; F=36000; G=900000},
    @{Row=31; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="2111"; F=36000; G=900000},
    @{Row=32; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="2110"; F=36000; G=900000},
    @{Row=33; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="2109"; F=36000; G=900000},
    @{Row=34; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="2108"; F=36000; G=900000},
    @{Row=35; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="2107"; F=36000; G=900000},
    @{Row=36; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="2106"; F=36000; G=900000},
    @{Row=37; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="2105"; F=36000; G=900000},
    @{Row=38; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="2104"; F=36000; G=900000},
    @{Row=39; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="2103"; F=36000; G=900000},
    @{Row=40; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="2102"; F=36000; G=900000},
    @{Row=41; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="2101"; F=36000; G=900000},
    @{Row=42; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="2012"; F=36000; G=900000},
    @{Row=43; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="2011"; F=36000; G=900000},
    @{Row=44; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="2010"; F=36000; G=900000},
    @{Row=45; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="2009"; F=36000; G=900000},
    @{Row=46; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="2008"; F=36000; G=900000},
    @{Row=47; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="2007"; F=36000; G=900000},
    @{Row=48; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="2006"; F=36000; G=900000},
    @{Row=49; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="2005"; F=36000; G=900000},
    @{Row=50; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="2004"; F=36000; G=900000},
    @{Row=51; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="2003"; F=36000; G=900000},
    @{Row=52; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="2002"; F=36000; G=900000},
    @{Row=53; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="2001"; F=36000; G=900000},
    @{Row=54; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="1912"; F=36000; G=900000},
    @{Row=55; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="1911"; F=36000; G=900000},
    @{Row=56; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="1910"; F=36000; G=900000},
    @{Row=57; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="1909"; F=36000; G=900000},
    @{Row=58; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="1908"; F=36000; G=900000},
    @{Row=59; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="1907"; F=36000; G=900000},
    @{Row=60; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="1906"; F=36000; G=900000},
    @{Row=61; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="1905"; F=36000; G=900000},
    @{Row=62; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="1904"; F=36000; G=900000},
    @{Row=63; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="1903"; F=36000; G=900000},
    @{Row=64; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="1902"; F=36000; G=900000},
    @{Row=65; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="1901"; F=36000; G=900000},
    @{Row=66; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="1812"; F=36000; G=900000},
    @{Row=67; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="1811"; F=36000; G=900000},
    @{Row=68; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="1810"; F=36000; G=900000},
    @{Row=69; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="1809"; F=36000; G=900000},
    @{Row=70; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="1808"; F=36000; G=900000},
    @{Row=71; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="1807"; F=36000; G=900000},
    @{Row=72; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="1806"; F=36000; G=900000},
    @{Row=73; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="1805"; F=36000; G=900000},
    @{Row=74; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="1804"; F=36000; G=900000},
    @{Row=75; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="1803"; F=36000; G=900000},
    @{Row=76; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="1802"; F=36000; G=900000},
    @{Row=77; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="1801"; F=36000; G=900000},
    @{Row=78; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="1712"; F=36000; G=900000},
    @{Row=79; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="1711"; F=36000; G=900000},
    @{Row=80; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="1710"; F=36000; G=900000},
    @{Row=81; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="1709"; F=36000; G=900000},
    @{Row=82; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="1708"; F=36000; G=900000},
    @{Row=83; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="1707"; F=36000; G=900000},
    @{Row=84; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="1706"; F=36000; G=900000},
    @{Row=85; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="1705"; F=36000; G=900000},
    @{Row=86; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="1704"; F=36000; G=900000},
    @{Row=87; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="1703"; F=36000; G=900000},
    @{Row=88; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="1702"; F=36000; G=900000},
    @{Row=89; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="1701"; F=36000; G=900000},
    @{Row=90; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="1612"; F=36000; G=900000},
    @{Row=91; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="1611"; F=36000; G=900000},
    @{Row=92; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="1610"; F=36000; G=900000},
    @{Row=93; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="1609"; F=36000; G=900000},
    @{Row=94; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="1608"; F=36000; G=900000},
    @{Row=95; C="73138679"; D="RICARDO ALVENIZ MADERO HERRERA"; E="1607"; F=36000; G=900000},
    @{Row=96; C="1045307877"; D="ADRIANA MILENA BARCASNEGRAS ORTIZ"; E="1909"; F=40000; G=1000000},
    @{Row=97; C="1045307877"; D="ADRIANA MILENA BARCASNEGRAS ORTIZ"; E="1908"; F=40000; G=1000000}
)

foreach ($row in $rows) {
    $r = $row.Row
    $ws.Cells.Item($r, 3).Value  = $row.C   # C -> N Doc Trabajador
    $ws.Cells.Item($r, 4).Value  = $row.D   # D -> Nombre Trabajador
    $ws.Cells.Item($r, 5).Value  = $row.E   # E -> Periodo Mora
    $ws.Cells.Item($r, 6).Value  = $row.F   # F -> Valor Mora
    $ws.Cells.Item($r, 7).Value  = $row.G   # G -> Salario Basico
}
